# Loan Origination with ML and RPA
# - Update the cached text of every "datetimeFigureOut" date field
#   (slide master + all slide layouts) from 14-02-2022 to 19-04-2022.
# - Rename/resize the "Mortgage or Loan Processing Process" textbox on
#   slide 3 to "Mortgage or Loan Processing Steps" and shrink its width
#   to match the new (shorter) caption.

$p = $ppt.ActivePresentation

function Set-DateFieldText {
    param($shapes, [string]$newText)

    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        $isDatePlaceholder = $false
        try {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                $isDatePlaceholder = $true
            }
        } catch {
            $isDatePlaceholder = $false
        }

        if ($isDatePlaceholder) {
            $shp.TextFrame.TextRange.Text = $newText
        }
    }
}

# Slide master.
$master = $p.SlideMaster
Set-DateFieldText $master.Shapes "19-04-2022"

# Every slide layout under the master.
$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    $layout = $layouts.Item($j)
    Set-DateFieldText $layout.Shapes "19-04-2022"
}

# Slide 3: rename + shrink the underline caption below the diagram.
$slide3 = $p.Slides.Item(3)
for ($i = 1; $i -le $slide3.Shapes.Count; $i++) {
    $shp = $slide3.Shapes.Item($i)
    if ($shp.Name -eq "TextBox 4") {
        $shp.TextFrame.TextRange.Text = "Mortgage or Loan Processing Steps"
        $shp.Width = 3553409 / 12700
    }
}
